$wb = $excel.ActiveWorkbook

# --- Service Contacts sheet: move "delivery_organisation_path" (currently
# --- the second-last column, R) so it sits right before "practitioner_key"
# --- (currently column D) -------------------------------------------------
$ws = $wb.Worksheets.Item("Service Contacts")

# Make room at D: shift practitioner_key..funding_source (D..Q) one column
# to the right.
$ws.Columns("D").Insert() | Out-Null

# delivery_organisation_path (header + its two data cells) is now at S
# (it shifted right along with everything else). Move it into the gap at D.
$ws.Range("S1:S3").Cut($ws.Range("D1:D3")) | Out-Null

# Remove the now-empty column it used to occupy.
$ws.Columns("S").Delete() | Out-Null

# Restore the column width (19 characters stored) on the relocated column,
# matching the width the source column used to carry.
$ws.Columns("D").ColumnWidth = 18.166666666666668

# Reflect the resulting selection / scroll position like a user who just
# performed this move would see: whole column D selected, view scrolled
# back to the left edge.
$ws.Application.Goto($ws.Range("A1"), $true) | Out-Null
$ws.Range("D:D").Select() | Out-Null

# --- K5 sheet: just update the remembered selection (no data changes) ----
$ws2 = $wb.Worksheets.Item("K5")
$ws2.Activate() | Out-Null
$ws2.Range("F1:F5").Select() | Out-Null
